$wb = $excel.ActiveWorkbook

# Update the "想去人数" (wishlist count) figures on the "展览" sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 353
$wsExhibit.Range("F4").Value = 4720

# Update the same figures on the aggregated "全部类型" sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 353
$wsAll.Range("F4").Value = 4720
